$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '82.290.21'
$ws.Range("E2").Value = '  +3.73%  '
$ws.Range("D3").Value = '3.199.88'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.292'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +24.15%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '3.188.60'
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.595'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000260'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +14.10%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.66%  '
$ws.Range("D15").Value = '3.769.52'
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.78'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '81.736.92'
$ws.Range("E17").Value = '  +3.23%  '
$ws.Range("D18").Value = '3.176.33'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.44%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '76.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '3.306.27'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.43%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000124'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '589.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.46%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.156'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("E37").Value = '  +16.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.16'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.58%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.409'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +21.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '159.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.25%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '188.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.78%  '
$ws.Range("E49").Value = '  +0.85%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.778'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.67%  '
